$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "matrices" score column (F) for every data row (2-13)
$ws.Range("F2").Value = 15.4707565882647
$ws.Range("F3").Value = 11.24716022619973
$ws.Range("F4").Value = 10.28396858577988
$ws.Range("F5").Value = 10.17922288895593
$ws.Range("F6").Value = 9.11564123155263
$ws.Range("F7").Value = 7.090080371342193
$ws.Range("F8").Value = 5.492612310008454
$ws.Range("F9").Value = 5.467819815978743
$ws.Range("F10").Value = 5.255963789825259
$ws.Range("F11").Value = 4.443507100185125
$ws.Range("F12").Value = 3.499128519851419
$ws.Range("F13").Value = 3.153071154622522

# Rows 4 and 5 swap places in the ranking (prolificid/name/gender move together)
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "60778ed0fde3e9c3a96f1d11"
$ws.Range("D4").Value = "Melissa"
$ws.Range("E4").Value = "female"

$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "60ba8ba51a5e0a105396888a"
$ws.Range("D5").Value = "Alfredo"
$ws.Range("E5").Value = "male"

# Rows 8 and 9 swap places in the ranking
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = "6024c18b094ac71dd93f4f5a"
$ws.Range("D8").Value = "Katherine"
$ws.Range("E8").Value = "female"

$ws.Range("B9").Value = 8
$ws.Range("C9").Value = "5f0142aa1eb1e528e7abce50"
$ws.Range("D9").Value = "Valeria"
$ws.Range("E9").Value = "female"
